$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2 = 0
    3 = 2
    4 = 2
    5 = 2
    6 = 0
    7 = 2
    8 = 4
    9 = 0
    10 = 0
    11 = 1
    12 = 4
    13 = 0
    14 = 2
    15 = 4
    16 = 3
    17 = 0
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 0
    23 = 3
    24 = 1
    25 = 1
    26 = 1
    27 = 1
    28 = 1
    29 = 0
    30 = 2
    31 = 0
    32 = 0
    33 = 0
    34 = 1
    35 = 1
    36 = 2
    37 = 2
    38 = 0
    39 = 1
    40 = 1
    41 = 1
    42 = 2
    43 = 2
    44 = 0
    45 = 1
    46 = 1
    47 = 3
    48 = 2
    49 = 0
    50 = 0
    51 = 1
    52 = 1
    53 = 0
    54 = 1
    55 = 0
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item([int]$row, 7).Value = $kValues[$row]
}

Write-Output "Done updating column G (K) for $($kValues.Count) rows"
